# Actualizacion 10 de Mayo
# Updates the statistics sheets with new figures and appends the list of
# "Rescatables" (students eligible for make-up/rescue exams) on sheet 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Estadisticos 1P"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D3").Value = 11
$ws1.Range("E3").Value = 1
$ws1.Range("H3").Value = 8.699999999999999

$ws1.Range("D4").Value = 7
$ws1.Range("E4").Value = 1
$ws1.Range("H4").Value = 7.8

# ---------------------------------------------------------------------
# Sheet "Estadisticos 2P"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 31
$ws2.Range("E2").Value = 13
$ws2.Range("F2").Value = 5
$ws2.Range("G2").Value = 13.89
$ws2.Range("H2").Value = 9

$ws2.Range("D3").Value = 18
$ws2.Range("E3").Value = 7
$ws2.Range("F3").Value = 13
$ws2.Range("G3").Value = 41.94
$ws2.Range("H3").Value = 8.699999999999999

$ws2.Range("D4").Value = 12
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Value = 9
$ws2.Range("G4").Value = 42.86
$ws2.Range("H4").Value = 8.1

# ---------------------------------------------------------------------
# Sheet "Estadisticos Final"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("H2").Value = 8

$ws3.Range("D3").Value = 11
$ws3.Range("F3").Value = 20
$ws3.Range("G3").Value = 64.52
$ws3.Range("H3").Value = 9

$ws3.Range("D4").Value = 7
$ws3.Range("F4").Value = 14
$ws3.Range("G4").Value = 66.67

# ---------------------------------------------------------------------
# Sheet "Rescatables" - append the list of students with pending exams
# (filled column by column, mirroring how the source data was pasted in)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$nc       = @(20330051920075, 20330051920062, 20330051920265, 20330051920268, 20330051920274, 20330051920278, 20330051920281)
$paterno  = @("CONTRERAS", "ROMERO", "CORONA", "GONZALEZ", "MARTINEZ", "RAMOS", "TELLEZ")
$materno  = @("GARCIA", "REYES", "HERNANDEZ", "FLORES", "RODRIGUEZ", "XOTLANIHUA", "OFICIAL")
$nombres  = @("JORGE HUMBERTO", "AMANDA MICHEL", "GUADALUPE", "JESUS HUMBERTO", "DANIEL ELEAZAR", "MARCO JOSAFAT", "MARISOL")
$materia  = @("LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II")
$grupo    = @("2AEV", "2AEV", "2APV", "2APV", "2APV", "2APV", "2APV")
$reprob   = @(2, 2, 2, 2, 2, 2, 2)

for ($i = 0; $i -lt $nc.Count; $i++) {
    $ws4.Cells.Item($i + 2, 1).Value = $nc[$i]
}
for ($i = 0; $i -lt $paterno.Count; $i++) {
    $ws4.Cells.Item($i + 2, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt $materno.Count; $i++) {
    $ws4.Cells.Item($i + 2, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt $nombres.Count; $i++) {
    $ws4.Cells.Item($i + 2, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt $materia.Count; $i++) {
    $ws4.Cells.Item($i + 2, 5).Value = $materia[$i]
}
for ($i = 0; $i -lt $grupo.Count; $i++) {
    $ws4.Cells.Item($i + 2, 6).Value = $grupo[$i]
}
for ($i = 0; $i -lt $reprob.Count; $i++) {
    $ws4.Cells.Item($i + 2, 7).Value = $reprob[$i]
}
